$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has 48 data rows (header rows 1-2, data rows 3-48).
# Add a new stage entry (#47) as row 49, matching the formatting of the
# previous data row (row 48).

$newRow = 49
$srcRow = 48

# Copy the fill/format from the previous row for columns A, C, D, E, F, G, H
# (column B keeps the un-filled "quote prefix" text style used for stageId
# values throughout the sheet, e.g. B3:B32).
foreach ($col in @(1,3,4,5,6,7,8)) {
    $ws.Cells.Item($srcRow, $col).Copy() | Out-Null
    $ws.Cells.Item($newRow, $col).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

# Populate the new stage's data.
$ws.Cells.Item($newRow, 1).Value = 47
$ws.Cells.Item($newRow, 2).Value = "'047"
$ws.Cells.Item($newRow, 3).Value = "MrHup"
$ws.Cells.Item($newRow, 4).Value = "인공 도시"
$ws.Cells.Item($newRow, 5).Value = "인공지능 기술로 그린 그림. 이제 인간이 설 자리는 어디인가?"
$ws.Cells.Item($newRow, 6).Value = 3
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = $false

# Leave selection on the last-edited cell, as in the authored workbook.
$ws.Cells.Item($newRow, 7).Select() | Out-Null
